$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 14:23"

# Apply the refreshed COVID-19 figures to each country currently holding that
# row (values are written in place; the subsequent Sort re-ranks every country
# by Casos totales, just like the live dashboard does on each data refresh).
# Row 4: Estados Unidos
$ws.Range("B4").Value = 4498887
$ws.Range("C4").Value = 544
$ws.Range("D4").Value = 2189592
$ws.Range("E4").Value = 2156937
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = 152358
# Row 6: India
$ws.Range("B6").Value = 1538899
$ws.Range("C6").Value = 6764
$ws.Range("D6").Value = 991205
$ws.Range("E6").Value = 513409
$ws.Range("G6").Value = 61
$ws.Range("H6").Value = 34285
# Row 26: Catar
$ws.Range("B26").Value = 110153
$ws.Range("C26").Value = 273
$ws.Range("D26").Value = 106849
$ws.Range("E26").Value = 3135
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 169
# Row 37: Bielorrusia
$ws.Range("B37").Value = 67518
$ws.Range("C37").Value = 152
$ws.Range("D37").Value = 61442
$ws.Range("E37").Value = 5528
$ws.Range("G37").Value = 5
$ws.Range("H37").Value = 548
# Row 40: Kuwait
$ws.Range("B40").Value = 65903
$ws.Range("C40").Value = 754
$ws.Range("D40").Value = 56467
$ws.Range("E40").Value = 8992
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 444
# Row 64: Uzbekistan
$ws.Range("E64").Value = 9775
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 129
# Row 68: Kenia
$ws.Range("B68").Value = 19125
$ws.Range("C68").Value = 544
$ws.Range("D68").Value = 8021
$ws.Range("E68").Value = 10793
$ws.Range("G68").Value = 12
$ws.Range("H68").Value = 311
# Row 73: Chequia
$ws.Range("B73").Value = 15827
$ws.Range("C73").Value = 28
$ws.Range("E73").Value = 4025
# Row 78: Dinamarca
$ws.Range("B78").Value = 13634
$ws.Range("C78").Value = 57
$ws.Range("D78").Value = 12485
$ws.Range("E78").Value = 535
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 614
# Row 82: Bosnia y Herzegovina
$ws.Range("B82").Value = 11127
$ws.Range("C82").Value = 361
$ws.Range("D82").Value = 5441
$ws.Range("E82").Value = 5370
$ws.Range("G82").Value = 19
$ws.Range("H82").Value = 316
# Row 84: Madagascar
$ws.Range("B84").Value = 10317
$ws.Range("C84").Value = 213
$ws.Range("D84").Value = 7117
$ws.Range("E84").Value = 3101
$ws.Range("G84").Value = 6
$ws.Range("H84").Value = 99
# Row 98: Zambia
$ws.Range("B98").Value = 5249
$ws.Range("C98").Value = 247
$ws.Range("D98").Value = 3285
$ws.Range("E98").Value = 1818
$ws.Range("G98").Value = 4
$ws.Range("H98").Value = 146
# Row 100: Croacia
$ws.Range("B100").Value = 4993
$ws.Range("C100").Value = 71
$ws.Range("D100").Value = 4099
$ws.Range("E100").Value = 753
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 141
# Row 144: Uganda
$ws.Range("B144").Value = 1140
$ws.Range("C144").Value = 5
$ws.Range("D144").Value = 1028
$ws.Range("E144").Value = 110
# Row 146: Burkina Faso
$ws.Range("D146").Value = 931
$ws.Range("E146").Value = 121
# Row 164: Vietnam
$ws.Range("B164").Value = 450
$ws.Range("C164").Value = 4
$ws.Range("D164").Value = 369
$ws.Range("E164").Value = 81

# Re-sort the data table (A3:H219, header in row 3) by Casos totales (column B) descending
$rng = $ws.Range("A3:H219")
$rng.Sort($ws.Range("B3"), 2, $null, $null, 1)
